# Fix "Errors in ATTRIBUTES REMOVED": the ATTRIBUTES column (D) mistakenly
# contained "REFILL_COUNT" at D6 -- remove it, shift the remaining
# attributes up one row, and append "CLASS DRUG" as the new final entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current ATTRIBUTES column values (D7:D13), skipping the
# erroneous "REFILL_COUNT" entry at D6, then shift them up by one row.
$attributes = @()
for ($r = 7; $r -le 13; $r++) {
    $attributes += $ws.Cells.Item($r, 4).Value2
}

for ($i = 0; $i -lt $attributes.Length; $i++) {
    $ws.Cells.Item(6 + $i, 4).Value = $attributes[$i]
}

# The new last row of the ATTRIBUTES column becomes "CLASS DRUG".
$ws.Range("D13").Value = "CLASS DRUG"

# Move the active selection to D13, matching the edited workbook.
$ws.Range("D13").Select()
